# Fix chapter header labels to include grade-level suffixes (e.g. "Geometry" -> "Geometry 7th")
# so that the three Term/Chapter names that are reused across grades 7, 8 and 9 are
# disambiguated, matching the JS side which already carries the grade suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-GradeHeader($cellAddress, $baseText, $grade, $fontName, $useBlackColor) {
    $full = "$baseText $grade" + "th"
    $r = $ws.Range($cellAddress)
    $r.Value = $full

    $boldLen = ($baseText + " " + $grade).Length
    $supStart = $boldLen + 1
    $supLen = 2

    $mainChars = $r.Characters(1, $boldLen)
    $mainChars.Font.Bold = $true
    $mainChars.Font.Size = 12
    $mainChars.Font.Name = $fontName
    if ($useBlackColor) {
        $mainChars.Font.Color = 0
    } else {
        $mainChars.Font.ColorIndex = -4105
    }

    $supChars = $r.Characters($supStart, $supLen)
    $supChars.Font.Bold = $true
    $supChars.Font.Size = 12
    $supChars.Font.Name = $fontName
    $supChars.Font.Superscript = $true
    if ($useBlackColor) {
        $supChars.Font.Color = 0
    } else {
        $supChars.Font.ColorIndex = -4105
    }
}

# Grade 7 chapters (Times New Roman)
Set-GradeHeader "B24" "Geometry" 7 "Times New Roman" $false
Set-GradeHeader "B32" "Practical Geometry" 7 "Times New Roman" $false

# Grade 8 chapters (Arial)
Set-GradeHeader "B38" "Algebra" 8 "Arial" $false
Set-GradeHeader "B62" "Practical Geometry" 8 "Arial" $false

# Grade 9 chapters (Arial)
Set-GradeHeader "B81" "Algebra" 9 "Arial" $false
Set-GradeHeader "B98" "Practical Geometry" 9 "Arial" $true
